$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to assign a text value to a cell while forcing text
# interpretation (prevents Excel from auto-converting numeric-looking
# strings like "92.00" into the number 92), then strips the temporary
# number-format override so the cell is left without an explicit style,
# matching the original workbook formatting.
function Set-TextValue($cellRef, [string]$val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

# Row 2
$ws.Range("D2").Value = '41.515.82'
$ws.Range("E2").Value = '  +0.05%  '

# Row 3
$ws.Range("D3").Value = '2.486.89'
$ws.Range("E3").Value = '  +0.94%  '

# Row 4
$ws.Range("E4").Value = '  +0.08%  '

# Row 5
Set-TextValue "D5" '310.58'
$ws.Range("E5").Value = '  -0.18%  '

# Row 6
Set-TextValue "D6" '92.00'
$ws.Range("E6").Value = '  -3.01%  '

# Row 7
Set-TextValue "D7" '0.537'
$ws.Range("E7").Value = '  -2.69%  '

# Row 8
$ws.Range("E8").Value = '  -0.13%  '

# Row 9
$ws.Range("E9").Value = '  -3.47%  '

# Row 10
Set-TextValue "D10" '32.09'
$ws.Range("E10").Value = '  -5.06%  '

# Row 11
$ws.Range("E11").Value = '  -1.03%  '

# Row 12
$ws.Range("E12").Value = '  +1.10%  '

# Row 13
$ws.Range("D13").Value = '2.868.84'
$ws.Range("E13").Value = '  +0.80%  '

# Row 14
Set-TextValue "D14" '6.78'
$ws.Range("E14").Value = '  -2.79%  '

# Row 15
$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '2.503.44'
$ws.Range("E15").Value = '  +1.82%  '

# Row 16
$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue "D16" '15.31'
$ws.Range("E16").Value = '  +4.95%  '

# Row 17
Set-TextValue "D17" '0.759'
$ws.Range("E17").Value = '  -3.88%  '

# Row 18
$ws.Range("D18").Value = '41.467.69'
$ws.Range("E18").Value = '  -0.25%  '

# Row 19
$ws.Range("E19").Value = '  -1.39%  '

# Row 20
$ws.Range("D20").Value = '0.0₃0917'
$ws.Range("E20").Value = '  -0.06%  '

# Row 21
Set-TextValue "D21" '70.35'
$ws.Range("E21").Value = '  +1.01%  '

# Row 22
Set-TextValue "D22" '11.07'
$ws.Range("E22").Value = '  -4.21%  '

# Row 23
Set-TextValue "D23" '234.25'
$ws.Range("E23").Value = '  -1.05%  '

# Row 24
$ws.Range("E24").Value = '  -3.37%  '

# Row 25
$ws.Range("E25").Value = '  -0.11%  '

# Row 26
$ws.Range("E26").Value = '  -2.79%  '

# Row 27
Set-TextValue "D27" '24.24'
$ws.Range("E27").Value = '  -1.95%  '

# Row 28
$ws.Range("E28").Value = '  +1.30%  '

# Row 29
Set-TextValue "D29" '9.58'
$ws.Range("E29").Value = '  -1.58%  '

# Row 30
Set-TextValue "D30" '36.24'
$ws.Range("E30").Value = '  -0.08%  '

# Row 31
Set-TextValue "D31" '153.44'
$ws.Range("E31").Value = '  +0.25%  '

# Row 32
$ws.Range("E32").Value = '  -4.70%  '

# Row 33
$ws.Range("B33").Value = 'Celestia'
$ws.Range("C33").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-TextValue "D33" '18.21'
$ws.Range("E33").Value = '  +6.11%  '

# Row 34
$ws.Range("B34").Value = 'WEMIXToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue "D34" '2.56'
$ws.Range("E34").Value = '  -2.26%  '

# Row 35
Set-TextValue "D35" '0.0755'
$ws.Range("E35").Value = '  +0.06%  '

# Row 36
Set-TextValue "D36" '2.49'
$ws.Range("E36").Value = '  -2.47%  '

# Row 37
Set-TextValue "D37" '2.97'
$ws.Range("E37").Value = '  -1.36%  '

# Row 38
$ws.Range("E38").Value = '  -3.04%  '

# Row 39
$ws.Range("E39").Value = '  -1.67%  '

# Row 41
Set-TextValue "D41" '4.03'
$ws.Range("E41").Value = '  +0.41%  '

# Row 42
$ws.Range("E42").Value = '  +0.36%  '

# Row 43
Set-TextValue "D43" '19.57'
$ws.Range("E43").Value = '  -8.09%  '

# Row 44
$ws.Range("D44").Value = '1.946.86'
$ws.Range("E44").Value = '  -1.80%  '

# Row 45
$ws.Range("E45").Value = '  -2.01%  '

# Row 46
$ws.Range("E46").Value = '  -4.16%  '

# Row 47
Set-TextValue "D47" '8.72'
$ws.Range("E47").Value = '  +0.11%  '

# Row 48
$ws.Range("D48").Value = '2.713.95'
$ws.Range("E48").Value = '  +0.63%  '

# Row 49
Set-TextValue "D49" '95.63'
$ws.Range("E49").Value = '  -1.93%  '

# Row 50
$ws.Range("E50").Value = '  -3.48%  '

# Row 51
Set-TextValue "D51" '66.28'
$ws.Range("E51").Value = '  -4.80%  '
